# "Generate Report for Handback"
#
# The 1ae9c0f9-4b1d-4f9f-a308-5ea809b2996b.md file has now been handed
# back (its translation is in sync with en-US) for both the zh-cn and
# de-de locales. Update the localization-status report:
#   - Overview sheet: flip that file's per-locale status columns.
#   - Per-locale sheets: flip Status, fill in the "Latest Target File" /
#     "Latest Handback File" columns (with hyperlinks matching the
#     existing Source/Target file links) and stamp the handback time.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: row 2 is the 1ae9c0f9...md file; zh-cn (B2) and
# de-de (C2) both move from "Ready for handoff" to "Handed back".
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $handedBack
$overview.Range("C2").Value = $handedBack

# ---------------------------------------------------------------------
# Helper that updates one locale sheet (zh-cn / de-de) for the
# 1ae9c0f9...md row (row 2): status, target/handback file + links,
# and the handback timestamp.
# ---------------------------------------------------------------------
function Update-LocaleSheet {
    param(
        [string]$SheetName,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Existing hyperlinks on row 2 (A2 = source .md, D2 = handoff .xlf)
    # give us the exact addresses/display text to reuse for the new
    # "Latest Target File" (F2) / "Latest Handback File" (G2) links.
    # (indexing the Hyperlinks collection directly is unreliable here,
    # so pull it into an array and match by range address instead.)
    $links = @($ws.Hyperlinks)
    $sourceAddress = $null
    $sourceDisplay = $null
    $xlfAddress = $null
    $xlfDisplay = $null
    foreach ($h in $links) {
        $addr = $h.Range.Address()
        if ($addr -eq '$A$2') {
            $sourceAddress = $h.Address
            $sourceDisplay = $h.TextToDisplay
        }
        if ($addr -eq '$D$2') {
            $xlfAddress = $h.Address
            $xlfDisplay = $h.TextToDisplay
        }
    }

    # Status flips to "Handed back: in sync with en-US"
    $ws.Range("C2").Value = $handedBack

    # Latest Target File (F2) -> same file/link as the source file (A2)
    $ws.Hyperlinks.Add($ws.Range("F2"), $sourceAddress, "", "", $sourceDisplay) | Out-Null

    # Latest Handback File (G2) -> same file/link as the handoff xlf (D2)
    $ws.Hyperlinks.Add($ws.Range("G2"), $xlfAddress, "", "", $xlfDisplay) | Out-Null

    # Latest Handback DateTime (H2)
    $ws.Range("H2").Value = $HandbackDateTime
}

Update-LocaleSheet "zh-cn" "2016-03-13 16:41:36"
Update-LocaleSheet "de-de" "2016-03-13 16:41:43"
